$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the measured length (D7) to the new value; it becomes a plain
# literal (no longer a formula) and dependent cells (D8:D11) recalc.
$ws.Range("D7").Value = 3.55

# The number format used by the D column values tightened from 3 to 4
# decimal places (0.000 -> 0.0000).
$ws.Range("D3:D11").NumberFormat = "0.0000"

# Move the active selection to D10 (cursor position when saved).
$ws.Range("D10").Select() | Out-Null
